$d = $word.ActiveDocument

# Locate the "Ver no Jupiter ..." paragraph (one of the trailing
# boilerplate/footer paragraphs that must be removed).
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "Ver no Jupiter Salvar em pdf Salvar em docx",
    $false, $false, $false, $false, $false,
    $true, 1, $false, "", 0)

if ($found) {
    $jupiterPara = $findRange.Paragraphs(1)

    # The paragraph immediately before it is the blank spacer paragraph,
    # and the paragraph immediately after it is the "(c) 2020 ..." credit
    # line. Together with the "Ver no Jupiter ..." paragraph itself these
    # three paragraphs (spacer, Jupiter line, copyright line) are removed,
    # while leaving the earlier "LOQ4219: ..." paragraph and the later
    # blank / page-break paragraphs untouched.
    $prevPara = $jupiterPara.Previous()
    $nextPara = $jupiterPara.Next()

    $deleteStart = $prevPara.Range.Start
    $deleteEnd = $nextPara.Range.End

    $d.Range($deleteStart, $deleteEnd).Delete()
}
